# Insert a new data row at row 64 (pushes existing rows 64-124 down to 65-125)
# and populate it with the new "Arveja Verde" observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("64:64").Insert()

$ws.Cells.Item(64, 1).Value  = 7
$ws.Cells.Item(64, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(64, 3).Value  = "Ñuble"
$ws.Cells.Item(64, 4).Value  = 45264
$ws.Cells.Item(64, 5).Value  = 16
$ws.Cells.Item(64, 6).Value  = 100112022
$ws.Cells.Item(64, 7).Value  = "Arveja Verde"
$ws.Cells.Item(64, 8).Value  = "Sin especificar"
$ws.Cells.Item(64, 9).Value  = "Primera"
$ws.Cells.Item(64, 10).Value = 50
$ws.Cells.Item(64, 11).Value = 24000
$ws.Cells.Item(64, 12).Value = 24000
$ws.Cells.Item(64, 13).Value = 24000
$ws.Cells.Item(64, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(64, 15).Value = "Región del Maule"
$ws.Cells.Item(64, 16).Value = 960
$ws.Cells.Item(64, 17).Value = 25
$ws.Cells.Item(64, 18).Value = "Hortaliza"
